# cmd/roi: add a "thumbnail" column (with per-shot thumbnail paths) to the
# roi test workbook, between the existing "duration" and "tags" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J ("duration" is I, "tags" was J -> becomes K),
# shifting the existing "tags" column (and everything to its right) one
# column to the right. Column formatting/styles of the existing cells are
# carried along automatically by Excel's column insert.
$ws.Columns.Item(10).Insert()

# New header cell for the inserted column.
$ws.Range("J1").Value = "thumbnail"

# Per-row thumbnail paths for the 10 data rows (rows 2-11).
$thumbs = @(
    "/vfx/thumbnail1.jpg",
    "/vfx/thumbnail2.jpg",
    "/vfx/thumbnail3.jpg",
    "/vfx/thumbnail4.jpg",
    "/vfx/thumbnail5.jpg",
    "/vfx/thumbnail6.jpg",
    "/vfx/thumbnail7.jpg",
    "/vfx/thumbnail8.jpg",
    "/vfx/thumbnail9.jpg",
    "/vfx/thumbnail10.jpg"
)

for ($i = 0; $i -lt $thumbs.Length; $i++) {
    $row = $i + 2
    $ws.Range("J$row").Value = $thumbs[$i]
}

# Style the new thumbnail data cells: italic "맑은 고딕" 11pt black, top
# vertical alignment (matches the rest of the data rows' look & feel).
$rng = $ws.Range("J2:J11")
$rng.Font.Italic = $true
$rng.Font.Name = "맑은 고딕"
$rng.Font.Size = 11
$rng.Font.Color = 0
$rng.VerticalAlignment = -4160

# Match the column widths used for the new layout: the thumbnail column is
# wider (paths are long), while the (now shifted) tags column keeps the
# width the duration/thumbnail columns had before.
$ws.Columns.Item(10).ColumnWidth = 17.065714285714286
$ws.Columns.Item(11).ColumnWidth = 11.035714285714286

# Leave the selection where the editor ended up after inserting the column
# and typing the new data.
$ws.Range("I12").Select()
